$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (ExpiryDate 2025-09-26)
$ws.Range("D3").Value = 118855.5795010159
$ws.Range("E3").Value = 0.0265629152391114
$ws.Range("F3").Value = 0.160564668370075
$ws.Range("G3").Value = -0.9380236697000435
$ws.Range("H3").Value = 10.29863907328862

# Row 5 (ExpiryDate 2025-11-28)
$ws.Range("D5").Value = 120516.6307708451
$ws.Range("E5").Value = 0.006901952297568422
$ws.Range("F5").Value = 0.2018481043281773
$ws.Range("G5").Value = -0.658992213993604
$ws.Range("H5").Value = 7.827816034389761

# Row 6 (ExpiryDate 2025-12-26)
$ws.Range("D6").Value = 120952.5460553045
$ws.Range("E6").Value = -0.008096012746467781
$ws.Range("F6").Value = 0.2418200688328248
$ws.Range("G6").Value = -1.233104038590722
$ws.Range("H6").Value = 10.29029350694815

# Row 7 (ExpiryDate 2026-01-30)
$ws.Range("D7").Value = 122265.2935309192
$ws.Range("E7").Value = -0.0341771056267732
$ws.Range("F7").Value = 0.3473456404648596
$ws.Range("G7").Value = -1.880787186326392
$ws.Range("H7").Value = 9.654951688688373

# Row 8 (ExpiryDate 2026-05-29)
$ws.Range("D8").Value = 122606.5212370685
$ws.Range("E8").Value = -0.03053021936026282
$ws.Range("F8").Value = 0.208334693899705
$ws.Range("G8").Value = -0.8598497897344676
$ws.Range("H8").Value = 6.965659856515654

# Row 9 (ExpiryDate 2026-06-26)
$ws.Range("D9").Value = 124174.0605675151
$ws.Range("E9").Value = -0.06696060077296213
$ws.Range("F9").Value = 0.3618027371121619
$ws.Range("G9").Value = -2.115727783351248
$ws.Range("H9").Value = 13.8010058011772

# Row 10 (ExpiryDate 2026-09-25)
$ws.Range("D10").Value = 125658.8905208434
$ws.Range("E10").Value = -0.1021339210113584
$ws.Range("F10").Value = 0.4341900863570095
$ws.Range("G10").Value = -1.898820265025255
$ws.Range("H10").Value = 9.725192333993247

# Row 11 (ExpiryDate 2026-12-24)
$ws.Range("D11").Value = 127728.4235471473
$ws.Range("E11").Value = -0.1788088683625924
$ws.Range("F11").Value = 0.761164008700197
$ws.Range("G11").Value = -2.611420253892948
$ws.Range("H11").Value = 12.75637441717517

# Row 17 (ExpiryDate 2025-08-15)
$ws.Range("D17").Value = 118151.9570776293
$ws.Range("E17").Value = 0.08976073939149479
$ws.Range("F17").Value = 0.1132811188622625
$ws.Range("G17").Value = -0.7824408353894894
$ws.Range("H17").Value = 6.530732845454605

# Row 18 (ExpiryDate 2025-08-22)
$ws.Range("D18").Value = 118205.3553134249
$ws.Range("E18").Value = 0.06127316797339297
$ws.Range("F18").Value = 0.1219528488774463
$ws.Range("G18").Value = -0.4946691996448183
$ws.Range("H18").Value = 7.238162325570007

# Row 19 (ExpiryDate 2025-09-05)
$ws.Range("D19").Value = 119005.9931964766
$ws.Range("E19").Value = 0.05704074090138368
$ws.Range("F19").Value = 0.1359210258614454
$ws.Range("G19").Value = -0.459548467173013
$ws.Range("H19").Value = 6.865036403677108
